# Insert a new "area" column (km^2 of each Polish voivodeship) as column B
# on every data sheet, shifting the existing co2/metan/n2o/so2/no/co columns
# one place to the right (B:G -> C:H).

$wb = $excel.ActiveWorkbook

# Area (powierzchnia) in km^2 per voivodeship - same values on every sheet,
# keyed by row number (row order is identical across all 5 sheets).
$areaByRow = @{
    2  = 19947   # dolnoslaskie
    3  = 17972   # kujawsko-pomorskie
    4  = 25122   # lubelskie
    5  = 13988   # lubuskie
    6  = 18219   # lodzkie
    7  = 15183   # malopolskie
    8  = 35558   # mazowieckie
    9  = 9412    # opolskie
    10 = 17846   # podkarpackie
    11 = 20187   # podlaskie
    12 = 18310   # pomorskie
    13 = 12333   # slaskie
    14 = 11711   # swietokrzyskie
    15 = 24173   # warminsko-mazurskie
    16 = 29826   # wielkopolskie
    17 = 22892   # zachodniopomorskie
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Insert a new blank column before column B - everything from B onward
    # (including styles) shifts right by one.
    $ws.Columns.Item(2).Insert()

    # Header
    $ws.Cells.Item(1, 2).Value = "area"

    # Data values - plain numbers, default ("Normalny") style, no inherited
    # borders/number format from the neighbouring cells.
    $dataRange = $ws.Range("B2:B17")
    $dataRange.Style = "Normalny"
    foreach ($r in $areaByRow.Keys) {
        $ws.Cells.Item($r, 2).Value = $areaByRow[$r]
    }

    # Selection: reflect that the new column B was just filled/selected.
    if ($i -eq 5) {
        $ws.Range("H8").Select()
    } else {
        $ws.Range("B1:B17").Select()
        $ws.Application.ActiveCell = $ws.Range("B1")
    }
}

$wb.Worksheets.Item(5).Activate()
